# Update "想去人数" (interest count) figures in column F across the
# 展览 / 演出 / 全部类型 sheets, matching the gh-pages data refresh at 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 5932   # was 5928
$ws.Range("F3").Value = 562    # was 561
$ws.Range("F4").Value = 1142   # was 1136
$ws.Range("F5").Value = 1077   # was 1073
$ws.Range("F6").Value = 859    # was 857
$ws.Range("F7").Value = 92     # was 91
$ws.Range("F8").Value = 54     # was 53
$ws.Range("F9").Value = 626    # was 624
$ws.Range("F10").Value = 69    # was 68
$ws.Range("F13").Value = 2114  # was 2100
$ws.Range("F14").Value = 1532  # was 1528
$ws.Range("F15").Value = 1182  # was 1172
$ws.Range("F17").Value = 217   # was 216
$ws.Range("F19").Value = 687   # was 679
$ws.Range("F20").Value = 248   # was 243
$ws.Range("F21").Value = 1076  # was 1075
$ws.Range("F24").Value = 3890  # was 3854
$ws.Range("F26").Value = 136   # was 135
$ws.Range("F28").Value = 174   # was 173
$ws.Range("F30").Value = 564   # was 560
$ws.Range("F35").Value = 338   # was 335
$ws.Range("F36").Value = 875   # was 872
$ws.Range("F38").Value = 78    # was 75
$ws.Range("F39").Value = 93    # was 92

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 762    # was 760

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5932   # was 5928
$ws.Range("F3").Value = 562    # was 561
$ws.Range("F4").Value = 1142   # was 1136
$ws.Range("F6").Value = 762    # was 760
$ws.Range("F7").Value = 1077   # was 1073
$ws.Range("F8").Value = 859    # was 857
$ws.Range("F11").Value = 92    # was 91
$ws.Range("F12").Value = 54    # was 53
$ws.Range("F13").Value = 626   # was 624
$ws.Range("F14").Value = 69    # was 68
$ws.Range("F18").Value = 2114  # was 2100
$ws.Range("F19").Value = 1532  # was 1528
$ws.Range("F20").Value = 1182  # was 1172
$ws.Range("F22").Value = 217   # was 216
$ws.Range("F25").Value = 687   # was 679
$ws.Range("F26").Value = 248   # was 243
$ws.Range("F27").Value = 1076  # was 1075
$ws.Range("F30").Value = 3890  # was 3855
$ws.Range("F32").Value = 136   # was 135
$ws.Range("F34").Value = 174   # was 173
$ws.Range("F36").Value = 564   # was 560
$ws.Range("F41").Value = 338   # was 335
$ws.Range("F42").Value = 875   # was 872
$ws.Range("F44").Value = 78    # was 75
$ws.Range("F45").Value = 93    # was 92
